$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71..171 down to 72..172
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new data point
$ws.Range("A71").Value = 9
$ws.Range("B71").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C71").Value = "Metropolitana"
$ws.Range("D71").Value = 44601
$ws.Range("E71").Value = 13
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100101
$ws.Range("H71").Value = "Berries"
$ws.Range("I71").Value = 100101001
$ws.Range("J71").Value = "Arándano (blue)"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 630
$ws.Range("N71").Value = 3500
$ws.Range("O71").Value = 4000
$ws.Range("P71").Value = 3722
$ws.Range("Q71").Value = "$/bandeja 2 kilos"
$ws.Range("R71").Value = "Región de O'Higgins"
$ws.Range("S71").Value = 1861
$ws.Range("T71").Value = 2
